$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "67.342.13"
Set-TextValue $ws.Range("E2") "  -1.11%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.620.35"
Set-TextValue $ws.Range("E3") "  -0.56%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.11%  "

# Row 5
Set-TextValue $ws.Range("D5") "589.79"
Set-TextValue $ws.Range("E5") "  +0.33%  "

# Row 6
Set-TextValue $ws.Range("D6") "183.99"
Set-TextValue $ws.Range("E6") "  +3.22%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.612"
Set-TextValue $ws.Range("E7") "  -2.06%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.26%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.672"
Set-TextValue $ws.Range("E9") "  -4.80%  "

# Row 10
Set-TextValue $ws.Range("E10") "  -9.22%  "

# Row 11
Set-TextValue $ws.Range("D11") "53.91"
Set-TextValue $ws.Range("E11") "  -3.07%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.0000253"
Set-TextValue $ws.Range("E12") "  -12.16%  "

# Row 13
Set-TextValue $ws.Range("E13") "  -5.21%  "

# Row 14
Set-TextValue $ws.Range("D14") "4.183.88"
Set-TextValue $ws.Range("E14") "  -1.29%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.611.55"
Set-TextValue $ws.Range("E15") "  -1.22%  "

# Row 16
Set-TextValue $ws.Range("E16") "  -0.07%  "

# Row 17
Set-TextValue $ws.Range("B17") "WrappedBTC"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "67.198.32"
Set-TextValue $ws.Range("E17") "  -1.01%  "

# Row 18
Set-TextValue $ws.Range("B18") "Chainlink"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D18") "18.39"
Set-TextValue $ws.Range("E18") "  -3.95%  "

# Row 19
Set-TextValue $ws.Range("D19") "12.22"
Set-TextValue $ws.Range("E19") "  -3.49%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -4.03%  "

# Row 21
Set-TextValue $ws.Range("D21") "392.82"
Set-TextValue $ws.Range("E21") "  -3.28%  "

# Row 22
Set-TextValue $ws.Range("D22") "4.33"
Set-TextValue $ws.Range("E22") "  -4.17%  "

# Row 23
Set-TextValue $ws.Range("D23") "84.88"
Set-TextValue $ws.Range("E23") "  -3.49%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -4.20%  "

# Row 25
Set-TextValue $ws.Range("D25") "12.27"
Set-TextValue $ws.Range("E25") "  -2.71%  "

# Row 26
Set-TextValue $ws.Range("D26") "6.06"
Set-TextValue $ws.Range("E26") "  +0.51%  "

# Row 27
Set-TextValue $ws.Range("D27") "10.36"
Set-TextValue $ws.Range("E27") "  -3.15%  "

# Row 28
Set-TextValue $ws.Range("D28") "3.60"
Set-TextValue $ws.Range("E28") "  -11.19%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.99"
Set-TextValue $ws.Range("E29") "  -4.19%  "

# Row 30
Set-TextValue $ws.Range("D30") "31.14"
Set-TextValue $ws.Range("E30") "  -3.89%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.77"
Set-TextValue $ws.Range("E31") "  -5.26%  "

# Row 32
Set-TextValue $ws.Range("B32") "Cosmos"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "11.91"
Set-TextValue $ws.Range("E32") "  -2.74%  "

# Row 33
Set-TextValue $ws.Range("B33") "OKB"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D33") "65.19"
Set-TextValue $ws.Range("E33") "  +1.39%  "

# Row 34
Set-TextValue $ws.Range("D34") "598.11"
Set-TextValue $ws.Range("E34") "  +1.44%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -2.81%  "

# Row 36
Set-TextValue $ws.Range("D36") "41.40"
Set-TextValue $ws.Range("E36") "  -3.16%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.24%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.999"
Set-TextValue $ws.Range("E38") "  -0.15%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.376"
Set-TextValue $ws.Range("E39") "  -4.86%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0₃0742"
Set-TextValue $ws.Range("E40") "  -15.45%  "

# Row 41
Set-TextValue $ws.Range("E41") "  -2.22%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -7.09%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.0412"
Set-TextValue $ws.Range("E43") "  -5.07%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.42"
Set-TextValue $ws.Range("E44") "  -9.71%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.691.77"
Set-TextValue $ws.Range("E45") "  +0.21%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -2.52%  "

# Row 47
Set-TextValue $ws.Range("D47") "3.02"
Set-TextValue $ws.Range("E47") "  -3.55%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.54"
Set-TextValue $ws.Range("E48") "  -5.85%  "

# Row 49
Set-TextValue $ws.Range("D49") "136.71"
Set-TextValue $ws.Range("E49") "  -2.56%  "

# Row 50
Set-TextValue $ws.Range("D50") "8.26"
Set-TextValue $ws.Range("E50") "  -7.46%  "

# Row 51
Set-TextValue $ws.Range("E51") "  -5.77%  "
